$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 18519974
$ws.Range("I62").Value = 30304414
$ws.Range("J62").Value = 1569.8572
$ws.Range("K62").Value = 30304414
$ws.Range("L62").Value = 1569.8572
$ws.Range("M62").Value = -30303790
$ws.Range("N62").Value = -2817.8572
$ws.Range("H65").Value = 18519974
$ws.Range("I65").Value = 30304414
$ws.Range("J65").Value = 1569.8572
$ws.Range("K65").Value = 151522070
$ws.Range("L65").Value = 7849.286
$ws.Range("M65").Value = -151518950
$ws.Range("N65").Value = -14089.286
$ws.Range("H76").Value = 3149777.2
$ws.Range("I76").Value = 4070294.8
$ws.Range("J76").Value = 4676.5835
$ws.Range("K76").Value = 4070294.8
$ws.Range("L76").Value = 4676.5835
$ws.Range("M76").Value = -4069979.8
$ws.Range("N76").Value = -5306.5835
$ws.Range("H79").Value = 3149777.2
$ws.Range("I79").Value = 4070294.8
$ws.Range("J79").Value = 4676.5835
$ws.Range("K79").Value = 4070294.8
$ws.Range("L79").Value = 4676.5835
$ws.Range("M79").Value = -4069202.8
$ws.Range("N79").Value = -6860.5835
$ws.Range("H106").Value = 45979716
$ws.Range("J106").Value = 66669468
$ws.Range("L106").Value = 66669468
$ws.Range("N106").Value = -66670730
$ws.Range("H112").Value = 1100.9474
$ws.Range("J112").Value = 1136.3529
$ws.Range("L112").Value = 3409.0587
$ws.Range("N112").Value = -5625.0587
$ws.Range("H116").Value = 6880.273
$ws.Range("I116").Value = 8832.857
$ws.Range("J116").Value = 3463.25
$ws.Range("K116").Value = 8832.857
$ws.Range("L116").Value = 3463.25
$ws.Range("M116").Value = -5390.857
$ws.Range("N116").Value = -10347.25
$ws.Range("H137").Value = 1287.3334
$ws.Range("I137").Value = 980.775
$ws.Range("J137").Value = 2402.0908
$ws.Range("K137").Value = 2942.325
$ws.Range("L137").Value = 7206.2724
$ws.Range("M137").Value = -392.3249999999998
$ws.Range("N137").Value = -12306.2724
$ws.Range("H138").Value = 2724.257
$ws.Range("I138").Value = 1164.4324
$ws.Range("K138").Value = 3493.2972
$ws.Range("M138").Value = 1646.7028

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5030.674
$ws.Range("I32").Value = 3545.405
$ws.Range("J32").Value = 12364.1875
$ws.Range("K32").Value = 3545.405
$ws.Range("L32").Value = 12364.1875
$ws.Range("M32").Value = -3258.405
$ws.Range("N32").Value = -12938.1875
$ws.Range("H45").Value = 11010.182
$ws.Range("I45").Value = 21102.4
$ws.Range("J45").Value = 2600
$ws.Range("K45").Value = 21102.4
$ws.Range("L45").Value = 2600
$ws.Range("M45").Value = -20725.4
$ws.Range("N45").Value = -3354
$ws.Range("H74").Value = 13514714
$ws.Range("I74").Value = 1021.1786
$ws.Range("J74").Value = 55557310
$ws.Range("K74").Value = 1021.1786
$ws.Range("L74").Value = 55557310
$ws.Range("M74").Value = -147.1786
$ws.Range("N74").Value = -55559058
$ws.Range("H77").Value = 13514714
$ws.Range("I77").Value = 1021.1786
$ws.Range("J77").Value = 55557310
$ws.Range("K77").Value = 5105.893
$ws.Range("L77").Value = 277786550
$ws.Range("M77").Value = -737.893
$ws.Range("N77").Value = -277795286
$ws.Range("H124").Value = 28235.285
$ws.Range("J124").Value = 28235.285
$ws.Range("L124").Value = 28235.285
$ws.Range("N124").Value = -38055.285
$ws.Range("H125").Value = 54147.2
$ws.Range("J125").Value = 54147.2
$ws.Range("L125").Value = 54147.2
$ws.Range("N125").Value = -63987.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 76924640
$ws.Range("I107").Value = 100001640
$ws.Range("J107").Value = 1300
$ws.Range("K107").Value = 100001640
$ws.Range("L107").Value = 1300
$ws.Range("M107").Value = -99999720
$ws.Range("N107").Value = -5140

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6875634.5
$ws.Range("I31").Value = 1776.7
$ws.Range("J31").Value = 14731472
$ws.Range("K31").Value = 1776.7
$ws.Range("L31").Value = 14731472
$ws.Range("M31").Value = -1481.7
$ws.Range("N31").Value = -14732062
$ws.Range("H34").Value = 6875634.5
$ws.Range("I34").Value = 1776.7
$ws.Range("J34").Value = 14731472
$ws.Range("K34").Value = 1776.7
$ws.Range("L34").Value = 14731472
$ws.Range("M34").Value = -1574.7
$ws.Range("N34").Value = -14731876

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 7333.8
$ws.Range("I122").Value = 419.9
$ws.Range("J122").Value = 10790.75
$ws.Range("K122").Value = 3779.1
$ws.Range("L122").Value = 97116.75
$ws.Range("M122").Value = -1329.1
$ws.Range("N122").Value = -102016.75
$ws.Range("H131").Value = 1471331.9
$ws.Range("I131").Value = 3030752
$ws.Range("J131").Value = 1021.45715
$ws.Range("K131").Value = 9092256
$ws.Range("L131").Value = 3064.37145
$ws.Range("M131").Value = -9087216
$ws.Range("N131").Value = -13144.37145
$ws.Range("H140").Value = 2763.2
$ws.Range("I140").Value = 2763.2
$ws.Range("K140").Value = 8289.599999999999
$ws.Range("M140").Value = -3109.599999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5657.8423
$ws.Range("I70").Value = 5269.231
$ws.Range("J70").Value = 6499.8335
$ws.Range("K70").Value = 5269.231
$ws.Range("L70").Value = 6499.8335
$ws.Range("M70").Value = -4999.231
$ws.Range("N70").Value = -7039.8335
$ws.Range("H73").Value = 5657.8423
$ws.Range("I73").Value = 5269.231
$ws.Range("J73").Value = 6499.8335
$ws.Range("K73").Value = 5269.231
$ws.Range("L73").Value = 6499.8335
$ws.Range("M73").Value = -4333.231
$ws.Range("N73").Value = -8371.833500000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 1016
$ws.Range("I30").Value = 1016
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 1016
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -908
$ws.Range("N30").ClearContents()
$ws.Range("H127").Value = 74499.5
$ws.Range("J127").Value = 74499.5
$ws.Range("L127").Value = 74499.5
$ws.Range("N127").Value = -84419.5
$ws.Range("H132").Value = 16197039
$ws.Range("I132").Value = 18430528
$ws.Range("J132").Value = 4250
$ws.Range("K132").Value = 55291584
$ws.Range("L132").Value = 12750
$ws.Range("M132").Value = -55289054
$ws.Range("N132").Value = -17810

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1143.375
$ws.Range("I132").Value = 905
$ws.Range("J132").Value = 1434.7222
$ws.Range("K132").Value = 2715
$ws.Range("L132").Value = 4304.1666
$ws.Range("M132").Value = -185
$ws.Range("N132").Value = -9364.1666
